$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 2 data rows (row 2 and row 3) below the header (row 1):
#   Row2: 23.01.2026,8,...  / 20.01.2026,1,...  / Stalis Donata / - / 4TFB|JA1 / Jezyk angielski / (empty)
#   Row3: 23.01.2026,9,...  / 20.01.2026,2,...  / Stalis Donata / - / 4TFB|JA1 / Jezyk angielski / (empty)
#
# A new record (Konczynska Malgorzata) needs to be inserted as the new row 2, pushing the
# existing two "Stalis Donata" rows down to row 3 and row 4.
#
# To avoid Excel copying the header's cell style into a freshly inserted row, we build the
# new layout by writing cell values directly (bottom-up) instead of using Rows.Insert().

# Step 1: push the current row 3 ("23.01.2026, 9, ..." record) down into row 4.
$ws.Cells.Item(4,1).Value = $ws.Cells.Item(3,1).Value2
$ws.Cells.Item(4,2).Value = $ws.Cells.Item(3,2).Value2
$ws.Cells.Item(4,3).Value = $ws.Cells.Item(3,3).Value2
$ws.Cells.Item(4,4).Value = $ws.Cells.Item(3,4).Value2
$ws.Cells.Item(4,5).Value = $ws.Cells.Item(3,5).Value2
$ws.Cells.Item(4,6).Value = $ws.Cells.Item(3,6).Value2

# Step 2: push the current row 2 ("23.01.2026, 8, ..." record) down into row 3.
$ws.Cells.Item(3,1).Value = $ws.Cells.Item(2,1).Value2
$ws.Cells.Item(3,2).Value = $ws.Cells.Item(2,2).Value2
$ws.Cells.Item(3,3).Value = $ws.Cells.Item(2,3).Value2
$ws.Cells.Item(3,4).Value = $ws.Cells.Item(2,4).Value2
$ws.Cells.Item(3,5).Value = $ws.Cells.Item(2,5).Value2
$ws.Cells.Item(3,6).Value = $ws.Cells.Item(2,6).Value2

# Step 3: write the new record into row 2.
$ws.Cells.Item(2,1).Value = "19.01.2026, 10, 15:45-16:30, sala: 31"
$ws.Cells.Item(2,2).Value = "19.01.2026, 5, 11:25-12:10, sala: 31"
$ws.Cells.Item(2,3).Value = "Kończyńska Małgorzata"
$ws.Cells.Item(2,4).Value = "-"
$ws.Cells.Item(2,5).Value = "2CB"
$ws.Cells.Item(2,6).Value = "Rozwój kompetencji zawodowych - dekoracje w cukiernictwie"

# Column width adjustments (best effort; widths expressed in Excel "character" units -
# the runtime rounds to the nearest 1/6 of a character, so these are the closest
# achievable values to the target stored widths of 31.71, 21.57 and 55.86).
# Column B is left untouched: splitting the old A:B shared column range by only
# resizing column A makes column B keep its original (already correct) width.
$ws.Columns.Item(1).ColumnWidth = 30.76
$ws.Columns.Item(3).ColumnWidth = 20.59
$ws.Columns.Item(6).ColumnWidth = 54.92

Write-Host "done"
